# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column, matching the header style of the existing headers, then
# fill in the per-row data values for rows 2-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) from the
# existing "IP" header cell (H1) onto the two new header cells so they
# match the look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2-13) ---
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 4

$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 7

$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 6

$ws.Range("I11").Value = 6
$ws.Range("J11").Value = 6

$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 6

$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 4
